$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "GROUP MEMBER": insert a new column A ("UserName") in front
# of the existing Date Invited/Accepted/Rejected/Left columns so the
# member rows line up with the USER sheet, same as the other
# per-user tables already do.
# ------------------------------------------------------------------
$wsGM = $wb.Worksheets.Item("GROUP MEMBER")
$wsGM.Columns("A").Insert()

$wsGM.Cells.Item(1,1).Value = "UserName"

$userNames = @("James","John","Robert","Michael","William","David","Richard","Joseph","Thomas","Charles","Christopher","Daniel","Matthew","Anthony","Donald","Mark","Paul","Steven","Andrew","Kenneth")
for ($i = 0; $i -lt $userNames.Count; $i++) {
    $wsGM.Cells.Item($i + 2, 1).Value = $userNames[$i]
}

# ------------------------------------------------------------------
# Sheet "ACTIVITY": ActivityId sequence in column A was renumbered
# (rows 3-10 bump up by one) and the leftover duplicate rows 11-21
# were cleared out, keeping only the date/time formatting in D:E.
# ------------------------------------------------------------------
$wsAct = $wb.Worksheets.Item("ACTIVITY")
$wsAct.Cells.Item(3,1).Value = 1092
$wsAct.Cells.Item(4,1).Value = 1093
$wsAct.Cells.Item(5,1).Value = 1094
$wsAct.Cells.Item(6,1).Value = 1095
$wsAct.Cells.Item(7,1).Value = 1096
$wsAct.Cells.Item(8,1).Value = 1097
$wsAct.Cells.Item(9,1).Value = 1098
$wsAct.Cells.Item(10,1).Value = 1099
$wsAct.Range("A11:E21").ClearContents()

# ------------------------------------------------------------------
# Sheet "CHALLENGE": the duplicated rows 8-21 were cleared out too,
# keeping only the date/time formatting in C:D.
# ------------------------------------------------------------------
$wsChal = $wb.Worksheets.Item("CHALLENGE")
$wsChal.Range("A8:B21").ClearContents()
